$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from E1 (header style) to F1, then set value
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate time_taken column F for rows 2-33
$ws.Range("F2").Value = "2021-10-05 10:52:10.233552"
$ws.Range("F3").Value = "2021-10-05 10:52:10.233566"
$ws.Range("F4").Value = "2021-10-05 10:52:10.233570"
$ws.Range("F5").Value = "2021-10-05 10:52:10.233573"
$ws.Range("F6").Value = "2021-10-05 10:52:10.233577"
$ws.Range("F7").Value = "2021-10-05 10:52:10.233580"
$ws.Range("F8").Value = "2021-10-05 10:52:10.233583"
$ws.Range("F9").Value = "2021-10-05 10:52:10.233586"
$ws.Range("F10").Value = "2021-10-05 10:52:10.233590"
$ws.Range("F11").Value = "2021-10-05 10:52:10.233593"
$ws.Range("F12").Value = "2021-10-05 10:52:10.233596"
$ws.Range("F13").Value = "2021-10-05 10:52:10.233599"
$ws.Range("F14").Value = "2021-10-05 10:52:10.233602"
$ws.Range("F15").Value = "2021-10-05 10:52:10.233606"
$ws.Range("F16").Value = "2021-10-05 10:52:10.233609"
$ws.Range("F17").Value = "2021-10-05 10:52:10.233612"
$ws.Range("F18").Value = "2021-10-05 10:52:10.233615"
$ws.Range("F19").Value = "2021-10-05 10:52:10.233618"
$ws.Range("F20").Value = "2021-10-05 10:52:10.233622"
$ws.Range("F21").Value = "2021-10-05 10:52:10.233625"
$ws.Range("F22").Value = "2021-10-05 10:52:10.233628"
$ws.Range("F23").Value = "2021-10-05 10:52:10.233631"
$ws.Range("F24").Value = "2021-10-05 10:52:10.233634"
$ws.Range("F25").Value = "2021-10-05 10:52:10.233637"
$ws.Range("F26").Value = "2021-10-05 10:52:10.233640"
$ws.Range("F27").Value = "2021-10-05 10:52:10.233643"
$ws.Range("F28").Value = "2021-10-05 10:52:10.233647"
$ws.Range("F29").Value = "2021-10-05 10:52:10.233650"
$ws.Range("F30").Value = "2021-10-05 10:52:10.233653"
$ws.Range("F31").Value = "2021-10-05 10:52:10.233656"
$ws.Range("F32").Value = "2021-10-05 10:52:10.233659"
$ws.Range("F33").Value = "2021-10-05 10:52:10.233663"

$excel.CutCopyMode = 0
